# Applies the "Microsoft Azure" skills-line tweak plus a batch of
# paragraph re-types that (as a side effect of Word re-flowing the runs)
# drop the now-unnecessary proofing-error markers (spellStart/spellEnd,
# gramStart/gramEnd) that were wrapping correctly-spelled proper nouns.

$d = $word.ActiveDocument
$apos = [char]0x2019

function Get-ParaContaining($doc, $marker) {
    foreach ($para in $doc.Paragraphs) {
        if ($para.Range.Text -like "*$marker*") {
            return $para
        }
    }
    return $null
}

function Retype-Paragraph($doc, $marker, $newText) {
    # Re-type the whole paragraph body so Word collapses it back down to a
    # single run and clears any stale proofing-error bookmarks that were
    # anchored to the old run boundaries.
    $para = Get-ParaContaining $doc $marker
    $rng = $para.Range
    $rng.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark
    $ok = $rng.Find.Execute($rng.Text, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    return $ok
}

# 1) Skills line: insert "Microsoft " right before the existing "Azure" run,
#    keeping every other run untouched (pure insertion, no re-typing).
$skills = Get-ParaContaining $d "Machine Learning, MongoDB"
$rng = $skills.Range
$rng.Find.Execute("Azure", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$d.TrackRevisions = $true
$rng.InsertBefore("Microsoft ")
$d.TrackRevisions = $false
$d.Revisions.Item(1).Accept()

# 2) Paragraphs whose proofing-error wrapped words get re-typed as-is so the
#    surrounding runs fold back into one contiguous run without the
#    <w:proofErr/> markers.
Retype-Paragraph $d "Palvogyi" ("Implemented the algorithm by Dang, Qi and Ye (2012), the algorithm by Fearnley, Palvogyi and Savani (2021), and a basic iteration algorithm to find Tarski" + $apos + "s fixed point in a complete lattice. ") | Out-Null

Retype-Paragraph $d "several experiments have been performed" "The algorithms were implemented in Python and several experiments have been performed to investigate the advantages of the algorithms in different scenarios. " | Out-Null

Retype-Paragraph $d "Pernoud" "Implemented the model by Eisenberg and Noe (2001) and the model by Jackson and Pernoud (2019) to compute clearing payments in financial networks. " | Out-Null

Retype-Paragraph $d "CartPole" "Successfully solved CartPole and Acrobot of OpenAI Gymnasium using Deep Q-networks and REINFROCE algorithm. " | Out-Null

Retype-Paragraph $d "Implemented using Python and" "Implemented using Python and PyTorch." | Out-Null

Retype-Paragraph $d "real-time HAR app" "This project aims to develop a real-time HAR app on Android devices using two sensors (Respeck and Thingy). " | Out-Null

Retype-Paragraph $d "Models used include" "Models used include Na$([char]0x00EF)ve Bayes, Maximum Entropy, Decision Tree, Random Forest, XGBoost, SVM, Multi-layer Perceptron, Recurrent Neural Network and Convolutional Neural Network." | Out-Null

Retype-Paragraph $d "smart contract written in Solidity" "Implemented as a smart contract written in Solidity and deployed on Ethereum testnet." | Out-Null

Retype-Paragraph $d "The robot was developed using" "The robot was developed using Webots and Python." | Out-Null

Write-Host "done"
